$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The shared status string "Ready for handoff" becomes "Handback transform failed"
# everywhere it is used (Overview!E3, Overview!F3, zh-cn!C3, de-de!C3).
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# New "Error Detail" values (column P) produced by the handback transform run.
$zhcn.Range("P3").Value = "Handback file name: bvzo1jhd.vxb is different with handoff file name: d9a05b3f-2d7b-4f3e-a406-2335368d3b54.a9356e99470b1b0bb8ee6d8aef9db82c9299a581.zh-cn."
$dede.Range("P3").Value = "Handback file name: bvzo1jhd.vxb is different with handoff file name: d9a05b3f-2d7b-4f3e-a406-2335368d3b54.a9356e99470b1b0bb8ee6d8aef9db82c9299a581.de-de."

# Widen the "Error Detail" column so the new messages are readable.
# (Excel stores column width with a +5/6 character padding offset, so to
# land exactly on a stored width of 40 we request 39 + 1/6.)
$newErrorDetailWidth = 39 + (1/6)
$zhcn.Columns.Item(16).ColumnWidth = $newErrorDetailWidth
$dede.Columns.Item(16).ColumnWidth = $newErrorDetailWidth
